$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.516.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.13"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "384.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.39"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  -1.26%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.589.46"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.097.30"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.992"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.92"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.553.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.26"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.63%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.29%  "

$ws.Range("E30").Value = "  -2.99%  "

$ws.Range("E31").Value = "  -2.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.36"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0480"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.15"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.06%  "

$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.07"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.14"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.066.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.946"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +20.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0325"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.00%  "
